$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a year-over-year table with one column per year (D=2007 ...
# P=2019). This adds the next column, Q, for year 2020, filling in its
# header (row 3) and all the data rows below it (rows 4-37). Row 34 is a
# blank spacer row in every other year column too, so Q34 is intentionally
# left without a value - it only needs the same formatting as P34.
$q2020 = [ordered]@{
  3  = 2020
  4  = 0.1
  5  = 0.1
  6  = 0.1
  7  = 0
  8  = 0
  9  = 0
  10 = 0
  11 = 0
  12 = 0
  13 = 0
  14 = 0.1
  15 = 0
  16 = 0
  17 = 0
  18 = 0
  19 = 0
  20 = 0
  21 = 0
  22 = 0
  23 = 0
  24 = 0
  25 = 0.1
  26 = 0.2
  27 = 0.1
  28 = 0.3
  29 = 0.4
  30 = 0.2
  31 = 0.2
  32 = 0.2
  33 = 0.1
  35 = 0
  36 = 0.1
  37 = 0.2
}

for ($row = 3; $row -le 37; $row++) {
  # Copy the existing P-column cell (same row) into Q so the new cell
  # picks up the identical number format / font / borders as the rest of
  # that row, then overwrite with this row's 2020 figure.
  $ws.Range("P$row").Copy($ws.Range("Q$row"))
  if ($q2020.Contains($row)) {
    $ws.Range("Q$row").Value = $q2020[$row]
  }
}
$excel.CutCopyMode = $false

# Leave the sheet's active selection on P30, matching the saved view state.
$ws.Range("P30").Select() | Out-Null
